$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the hyperlinks that were previously attached to C2:C4 (the
# register test data no longer mails these addresses).
$ws.Hyperlinks.Delete()

# Update the email values for the three "register" rows.
$ws.Range("C2").Value = "saurav.mehta1@testmail.com"
$ws.Range("C3").Value = "sofia.fernandez2@testmail.com"
$ws.Range("C4").Value = "natalie.khan2@stmail.com"

# The "Email" field (column F) now carries SQL-injection-style payloads
# used to exercise the registration form's input validation.
$ws.Range("F2").Value = "OR 1=1; -- , testuser@example.com"
$ws.Range("F3").Value = "admin' DROP TABLE users; -- , random.email123@mail.com"
$ws.Range("F4").Value = "{""payload"":""' OR 1=1; --"",""email"":""testuser@example.com""}"

$ws.Range("F2").WrapText = $false
$ws.Range("F3").WrapText = $false
$ws.Range("F4").WrapText = $true

$ws.Columns.Item(6).ColumnWidth = 50.77734375
